# Updated cryptos list: refresh Price (col D) and Volume/1h change (col E)
# figures scraped for the new run, cell by cell, mirroring the upstream
# generator which stores every D/E value as plain inline text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.185.90'
$ws.Range('E2').Value = '  -0.59%  '
$ws.Range('D3').Value = '2.298.35'
$ws.Range('E3').Value = '  -1.63%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '106.27'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.48%  '
$ws.Range('E7').Value = '  -1.47%  '
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.611'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.18'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0912'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.40'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.107'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.979'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.37'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.96%  '
$ws.Range('D16').Value = '2.649.26'
$ws.Range('E16').Value = '  -1.47%  '
$ws.Range('D17').Value = '2.304.75'
$ws.Range('E17').Value = '  -1.50%  '
$ws.Range('D18').Value = '42.030.72'
$ws.Range('E18').Value = '  -1.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.64'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.91'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.87%  '
$ws.Range('E22').Value = '  -0.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '260.38'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.66%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.32'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.81'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.68%  '
$ws.Range('E26').Value = '  +0.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.94%  '
$ws.Range('E28').Value = '  +2.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '22.80'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.50'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '166.28'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0893'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.53%  '
$ws.Range('E33').Value = '  -1.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.89'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.52%  '
$ws.Range('E35').Value = '  +7.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.130'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('E37').Value = '  +2.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.93'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +12.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0353'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.62'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.85%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '100.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +18.02%  '
$ws.Range('E42').Value = '  +1.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '71.16'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.86%  '
$ws.Range('E44').Value = '  -1.42%  '
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.28'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '113.35'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.55%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '78.81'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +8.64%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.17'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.33'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.65%  '
$ws.Range('E51').Value = '  +3.13%  '
